# Applies the JOURNAL.docx edit described by the diff:
#  - wraps a handful of "odd" tokens (sw.js, JQuery, Impromtu,
#    reconfirm_password(), set_bt_pwd()) in <w:proofErr> markers by
#    splitting their paragraph's single run into several runs
#  - appends one new, empty paragraph at the end of the document

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Replace-ParagraphRuns {
    # NOTE: this emulated PowerShell host only binds *positional*
    # arguments reliably for user-defined functions; -Name Value style
    # binding silently leaves the parameter empty. Always call this
    # positionally: Replace-ParagraphRuns <index> <xml>
    param(
        [int]$ParaIndex,
        [string]$InnerXml
    )
    $para = $d.Paragraphs($ParaIndex).Range
    # Exclude the trailing paragraph-mark character so the paragraph
    # itself (and its pPr) is preserved; only the run content inside
    # is replaced by InsertXML.
    $target = $d.Range($para.Start, $para.End - 1)
    $xml = "<w:p xmlns:w='$wNs'>" + $InnerXml + "</w:p>"
    $target.InsertXML($xml)
}

# --- Paragraph 2: "Added jura to the CSS files section and added to sw.js"
$inner = ""
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>Added jura to the CSS files section and added to </w:t></w:r>"
$inner += "<w:proofErr w:type='gramStart'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>sw.js</w:t></w:r>"
$inner += "<w:proofErr w:type='gramEnd'/>"
Replace-ParagraphRuns 2 $inner

# --- Paragraph 3: "Added JQuery and Impromtu plugins to introduce prompt with addons:"
$inner = ""
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>Added </w:t></w:r>"
$inner += "<w:proofErr w:type='spellStart'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>JQuery</w:t></w:r>"
$inner += "<w:proofErr w:type='spellEnd'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> and </w:t></w:r>"
$inner += "<w:proofErr w:type='spellStart'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Impromtu</w:t></w:r>"
$inner += "<w:proofErr w:type='spellEnd'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> plugins to introduce prompt with addons:</w:t></w:r>"
Replace-ParagraphRuns 3 $inner

# --- Paragraph 6: "Bluetooth password setting is nor controlled by reconfirm_password()."
$inner = ""
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>Bluetooth password setting is nor controlled by </w:t></w:r>"
$inner += "<w:proofErr w:type='spellStart'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>reconfirm_</w:t></w:r>"
$inner += "<w:proofErr w:type='gramStart'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>password</w:t></w:r>"
$inner += "<w:proofErr w:type='spellEnd'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>(</w:t></w:r>"
$inner += "<w:proofErr w:type='gramEnd'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>).</w:t></w:r>"
Replace-ParagraphRuns 6 $inner

# --- Paragraph 8: "If matched, the password set function is called (set_bt_pwd()). If failed, an error message is thrown."
$inner = ""
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>If matched, the password set function is called (</w:t></w:r>"
$inner += "<w:proofErr w:type='spellStart'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>set_bt_</w:t></w:r>"
$inner += "<w:proofErr w:type='gramStart'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>pwd</w:t></w:r>"
$inner += "<w:proofErr w:type='spellEnd'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>(</w:t></w:r>"
$inner += "<w:proofErr w:type='gramEnd'/>"
$inner += "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>)). If failed, an error message is thrown.</w:t></w:r>"
Replace-ParagraphRuns 8 $inner

# --- Append a new, empty paragraph after the last paragraph
#     ("Version updated to 1.0.8"), right before the sectPr.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex).Range
$lastPara.InsertParagraphAfter()
$newPara = $d.Paragraphs($lastIndex + 1).Range
$newXml = "<w:p xmlns:w='$wNs'><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr></w:p>"
$newPara.InsertXML($newXml)
